# Weekly CompStat crime-data refresh (5th Precinct, week of 3/6/2023-3/12/2023).
# Updates header (volume/date) text and the Week-to-Date / 28-Day / Year-to-Date
# / 2-Year crime-count + percent-change table (rows 14-30).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text (rich-text shared strings; in-place substring replace) ---
$ws.Range("A8").Replace("9", "10")
$ws.Range("C9").Replace("2/27/2023", "3/6/2023")
$ws.Range("C9").Replace("3/5/2023", "3/12/2023")

# --- Donor template cells (stable styles/content throughout this script) ---
$donorNum15 = $ws.Cells.Item(15, 6)   # F15: style 15 (integer count format)
$donorNum16 = $ws.Cells.Item(14, 8)   # H14: style 16 (percent-change format)
$donorPh20  = $ws.Cells.Item(14, 3)   # C14: style 14, shared text "0"
$donorPh21  = $ws.Cells.Item(14, 12)  # L14: style 14, shared text "***.*"

# --- Cells that change representation (number <-> text placeholder) ---
$donorNum15.Copy($ws.Cells.Item(14, 4))
$ws.Cells.Item(14, 4).Value = 1
$donorNum16.Copy($ws.Cells.Item(14, 5))
$ws.Cells.Item(14, 5).Value = -100
$donorPh20.Copy($ws.Cells.Item(15, 3))
$donorNum15.Copy($ws.Cells.Item(18, 4))
$ws.Cells.Item(18, 4).Value = 1
$donorNum16.Copy($ws.Cells.Item(18, 5))
$ws.Cells.Item(18, 5).Value = 200
$donorPh20.Copy($ws.Cells.Item(20, 3))
$donorPh20.Copy($ws.Cells.Item(20, 4))
$donorPh21.Copy($ws.Cells.Item(20, 5))
$donorPh20.Copy($ws.Cells.Item(22, 3))
$donorPh20.Copy($ws.Cells.Item(22, 4))
$donorPh21.Copy($ws.Cells.Item(22, 5))
$donorPh20.Copy($ws.Cells.Item(23, 3))
$donorNum15.Copy($ws.Cells.Item(23, 4))
$ws.Cells.Item(23, 4).Value = 1
$donorNum16.Copy($ws.Cells.Item(23, 5))
$ws.Cells.Item(23, 5).Value = -100
$donorPh20.Copy($ws.Cells.Item(26, 3))
$donorNum15.Copy($ws.Cells.Item(28, 4))
$ws.Cells.Item(28, 4).Value = 1
$donorNum16.Copy($ws.Cells.Item(28, 5))
$ws.Cells.Item(28, 5).Value = -100
$donorNum15.Copy($ws.Cells.Item(28, 7))
$ws.Cells.Item(28, 7).Value = 1
$donorNum16.Copy($ws.Cells.Item(28, 8))
$ws.Cells.Item(28, 8).Value = -100
$donorNum15.Copy($ws.Cells.Item(28, 10))
$ws.Cells.Item(28, 10).Value = 1
$donorNum16.Copy($ws.Cells.Item(28, 11))
$ws.Cells.Item(28, 11).Value = 100
$donorNum15.Copy($ws.Cells.Item(29, 4))
$ws.Cells.Item(29, 4).Value = 1
$donorNum16.Copy($ws.Cells.Item(29, 5))
$ws.Cells.Item(29, 5).Value = -100
$donorNum15.Copy($ws.Cells.Item(29, 7))
$ws.Cells.Item(29, 7).Value = 1
$donorNum16.Copy($ws.Cells.Item(29, 8))
$ws.Cells.Item(29, 8).Value = -100
$donorNum15.Copy($ws.Cells.Item(29, 10))
$ws.Cells.Item(29, 10).Value = 1
$donorNum16.Copy($ws.Cells.Item(29, 11))
$ws.Cells.Item(29, 11).Value = 100
$donorPh20.Copy($ws.Cells.Item(30, 4))
$donorPh21.Copy($ws.Cells.Item(30, 5))

# --- Plain numeric value updates (style/format unchanged) ---
$ws.Cells.Item(14, 7).Value = 2
$ws.Cells.Item(14, 10).Value = 2
$ws.Cells.Item(15, 12).Value = -33.333333333333
$ws.Cells.Item(16, 3).Value = 2
$ws.Cells.Item(16, 4).Value = 1
$ws.Cells.Item(16, 5).Value = 100
$ws.Cells.Item(16, 9).Value = 31
$ws.Cells.Item(16, 10).Value = 13
$ws.Cells.Item(16, 11).Value = 138.461538461538
$ws.Cells.Item(16, 12).Value = 106.666666666667
$ws.Cells.Item(16, 13).Value = 72.222222222222
$ws.Cells.Item(16, 14).Value = -75.2
$ws.Cells.Item(17, 3).Value = 2
$ws.Cells.Item(17, 4).Value = 11
$ws.Cells.Item(17, 5).Value = -81.818181818181
$ws.Cells.Item(17, 7).Value = 25
$ws.Cells.Item(17, 8).Value = -52
$ws.Cells.Item(17, 9).Value = 34
$ws.Cells.Item(17, 10).Value = 43
$ws.Cells.Item(17, 11).Value = -20.930232558139
$ws.Cells.Item(17, 12).Value = 126.666666666667
$ws.Cells.Item(17, 13).Value = 54.545454545454
$ws.Cells.Item(17, 14).Value = -20.930232558139
$ws.Cells.Item(18, 6).Value = 11
$ws.Cells.Item(18, 7).Value = 6
$ws.Cells.Item(18, 8).Value = 83.333333333333
$ws.Cells.Item(18, 9).Value = 27
$ws.Cells.Item(18, 10).Value = 27
$ws.Cells.Item(18, 11).Value = 0
$ws.Cells.Item(18, 12).Value = 107.692307692308
$ws.Cells.Item(18, 13).Value = -27.027027027027
$ws.Cells.Item(18, 14).Value = -82.911392405063
$ws.Cells.Item(19, 3).Value = 12
$ws.Cells.Item(19, 4).Value = 10
$ws.Cells.Item(19, 5).Value = 20
$ws.Cells.Item(19, 6).Value = 51
$ws.Cells.Item(19, 7).Value = 41
$ws.Cells.Item(19, 8).Value = 24.390243902439
$ws.Cells.Item(19, 9).Value = 128
$ws.Cells.Item(19, 10).Value = 96
$ws.Cells.Item(19, 11).Value = 33.333333333333
$ws.Cells.Item(19, 12).Value = 82.857142857142
$ws.Cells.Item(19, 13).Value = 28
$ws.Cells.Item(19, 14).Value = -51.515151515151
$ws.Cells.Item(20, 6).Value = 4
$ws.Cells.Item(20, 8).Value = 33.333333333333
$ws.Cells.Item(20, 14).Value = -86.842105263157
$ws.Cells.Item(21, 3).Value = 19
$ws.Cells.Item(21, 4).Value = 24
$ws.Cells.Item(21, 5).Value = -20.833333333333
$ws.Cells.Item(21, 6).Value = 93
$ws.Cells.Item(21, 7).Value = 86
$ws.Cells.Item(21, 8).Value = 8.13953488372
$ws.Cells.Item(21, 9).Value = 232
$ws.Cells.Item(21, 10).Value = 189
$ws.Cells.Item(21, 11).Value = 22.751322751322
$ws.Cells.Item(21, 12).Value = 91.735537190082
$ws.Cells.Item(21, 13).Value = 28.17679558011
$ws.Cells.Item(21, 14).Value = -65.321375186846
$ws.Cells.Item(22, 6).Value = 1
$ws.Cells.Item(22, 8).Value = -88.888888888888
$ws.Cells.Item(22, 12).Value = -62.5
$ws.Cells.Item(23, 6).Value = 1
$ws.Cells.Item(23, 7).Value = 6
$ws.Cells.Item(23, 8).Value = -83.333333333333
$ws.Cells.Item(23, 10).Value = 7
$ws.Cells.Item(23, 11).Value = 14.285714285714
$ws.Cells.Item(23, 13).Value = 14.285714285714
$ws.Cells.Item(24, 3).Value = 12
$ws.Cells.Item(24, 4).Value = 11
$ws.Cells.Item(24, 5).Value = 9.090909090909
$ws.Cells.Item(24, 6).Value = 75
$ws.Cells.Item(24, 7).Value = 53
$ws.Cells.Item(24, 8).Value = 41.509433962264
$ws.Cells.Item(24, 9).Value = 197
$ws.Cells.Item(24, 10).Value = 166
$ws.Cells.Item(24, 11).Value = 18.67469879518
$ws.Cells.Item(24, 12).Value = 109.574468085106
$ws.Cells.Item(24, 13).Value = 11.299435028248
$ws.Cells.Item(25, 3).Value = 6
$ws.Cells.Item(25, 5).Value = -50
$ws.Cells.Item(25, 6).Value = 20
$ws.Cells.Item(25, 7).Value = 39
$ws.Cells.Item(25, 8).Value = -48.717948717948
$ws.Cells.Item(25, 9).Value = 58
$ws.Cells.Item(25, 10).Value = 70
$ws.Cells.Item(25, 11).Value = -17.142857142857
$ws.Cells.Item(25, 12).Value = 41.463414634146
$ws.Cells.Item(25, 13).Value = 28.888888888888
$ws.Cells.Item(26, 12).Value = -40
$ws.Cells.Item(27, 4).Value = 2
$ws.Cells.Item(27, 5).Value = -50
$ws.Cells.Item(27, 6).Value = 4
$ws.Cells.Item(27, 8).Value = -55.555555555555
$ws.Cells.Item(27, 9).Value = 8
$ws.Cells.Item(27, 10).Value = 13
$ws.Cells.Item(27, 11).Value = -38.461538461538
$ws.Cells.Item(27, 12).Value = -20
$ws.Cells.Item(30, 7).Value = 3
